$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

$ws.Cells.Item(2, 9).Value = 'b'
$ws.Cells.Item(2, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(8, 9).Value = 'aa'
$ws.Cells.Item(8, 10).Value = 'Agree/Accept'
$ws.Cells.Item(29, 9).Value = 'aa'
$ws.Cells.Item(29, 10).Value = 'Agree/Accept'
$ws.Cells.Item(34, 9).Value = 'sv'
$ws.Cells.Item(34, 10).Value = 'Statement-opinion'
$ws.Cells.Item(37, 9).Value = 'sd'
$ws.Cells.Item(37, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(40, 9).Value = 'aa'
$ws.Cells.Item(40, 10).Value = 'Agree/Accept'
$ws.Cells.Item(60, 9).Value = 'sd'
$ws.Cells.Item(60, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(78, 9).Value = 'b'
$ws.Cells.Item(78, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(90, 9).Value = 'sv'
$ws.Cells.Item(90, 10).Value = 'Statement-opinion'
$ws.Cells.Item(102, 9).Value = 'sv'
$ws.Cells.Item(102, 10).Value = 'Statement-opinion'
$ws.Cells.Item(109, 9).Value = 'aa'
$ws.Cells.Item(109, 10).Value = 'Agree/Accept'
$ws.Cells.Item(111, 9).Value = 'aa'
$ws.Cells.Item(111, 10).Value = 'Agree/Accept'
$ws.Cells.Item(117, 9).Value = 'aa'
$ws.Cells.Item(117, 10).Value = 'Agree/Accept'
$ws.Cells.Item(118, 9).Value = 'b'
$ws.Cells.Item(118, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(120, 9).Value = 'ba'
$ws.Cells.Item(120, 10).Value = 'Appreciation'
$ws.Cells.Item(131, 9).Value = 'ba'
$ws.Cells.Item(131, 10).Value = 'Appreciation'
$ws.Cells.Item(135, 9).Value = 'sd'
$ws.Cells.Item(135, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(142, 9).Value = 'sv'
$ws.Cells.Item(142, 10).Value = 'Statement-opinion'
$ws.Cells.Item(154, 9).Value = 'sv'
$ws.Cells.Item(154, 10).Value = 'Statement-opinion'
$ws.Cells.Item(157, 9).Value = 'aa'
$ws.Cells.Item(157, 10).Value = 'Agree/Accept'
$ws.Cells.Item(188, 9).Value = 'sd'
$ws.Cells.Item(188, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(194, 9).Value = '%'
$ws.Cells.Item(194, 10).Value = 'Uninterpretable'
$ws.Cells.Item(199, 9).Value = 'b'
$ws.Cells.Item(199, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(202, 9).Value = 'sd'
$ws.Cells.Item(202, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(210, 9).Value = 'aa'
$ws.Cells.Item(210, 10).Value = 'Agree/Accept'
$ws.Cells.Item(217, 9).Value = 'sd'
$ws.Cells.Item(217, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(221, 9).Value = 'aa'
$ws.Cells.Item(221, 10).Value = 'Agree/Accept'
$ws.Cells.Item(223, 9).Value = 'b'
$ws.Cells.Item(223, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(232, 9).Value = 'b'
$ws.Cells.Item(232, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(241, 9).Value = 'b'
$ws.Cells.Item(241, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(259, 9).Value = 'sd'
$ws.Cells.Item(259, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(263, 9).Value = 'ba'
$ws.Cells.Item(263, 10).Value = 'Appreciation'
$ws.Cells.Item(276, 9).Value = 'aa'
$ws.Cells.Item(276, 10).Value = 'Agree/Accept'
$ws.Cells.Item(280, 9).Value = 'aa'
$ws.Cells.Item(280, 10).Value = 'Agree/Accept'
$ws.Cells.Item(289, 9).Value = 'aa'
$ws.Cells.Item(289, 10).Value = 'Agree/Accept'
$ws.Cells.Item(296, 9).Value = 'sd'
$ws.Cells.Item(296, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(303, 9).Value = 'b'
$ws.Cells.Item(303, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(314, 9).Value = 'sd'
$ws.Cells.Item(314, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(318, 9).Value = 'b'
$ws.Cells.Item(318, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(319, 9).Value = 'sd'
$ws.Cells.Item(319, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(322, 9).Value = 'ba'
$ws.Cells.Item(322, 10).Value = 'Appreciation'
$ws.Cells.Item(348, 9).Value = 'sv'
$ws.Cells.Item(348, 10).Value = 'Statement-opinion'
$ws.Cells.Item(349, 9).Value = 'sd'
$ws.Cells.Item(349, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(351, 9).Value = '%'
$ws.Cells.Item(351, 10).Value = 'Uninterpretable'
$ws.Cells.Item(352, 9).Value = 'sv'
$ws.Cells.Item(352, 10).Value = 'Statement-opinion'
$ws.Cells.Item(370, 9).Value = 'b'
$ws.Cells.Item(370, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(403, 9).Value = 'sv'
$ws.Cells.Item(403, 10).Value = 'Statement-opinion'
$ws.Cells.Item(408, 9).Value = 'b'
$ws.Cells.Item(408, 10).Value = 'Acknowledge (Backchannel)'